$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: "Pura zip-paketti johonkin hakemistoon, esim c:\kmean" was split
# across a run, a (now removed) "_GoBack" bookmark, and a trailing run "s".
# Delete the stray bookmark and reflow the text into one merged run reading
# "...esim c:\kmeans".
# -------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Content.Find.Execute("Pura zip-paketti johonkin hakemistoon, esim c:\kmeans", $true, $false, $false, $false, $false, $true, 1, $false, "Pura zip-paketti johonkin hakemistoon, esim c:\kmeans", 2)

# -------------------------------------------------------------------------
# Change 2: Add the NetBeans version number to the "Netbeans" heading, i.e.
# turn "Netbeans" into "Netbeans (8.x)" while keeping the existing
# _Toc257918838 bookmark wrapped tightly around just "Netbeans", and leave a
# fresh "_GoBack" bookmark collapsed at the very end of the paragraph
# (after " (8.x)").
# -------------------------------------------------------------------------

# Locate the "Netbeans" heading paragraph (the one using Heading2 style whose
# text is exactly "Netbeans").
$netbeansPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Netbeans") {
        $netbeansPara = $p
        break
    }
}

$headingStart = $netbeansPara.Range.Start

# Replace "Netbeans" with "Netbeans (8.x)" -- scoping the Find to this
# paragraph only keeps the run's existing character formatting (w:lang) and
# avoids touching the "Avaa Netbeans" bullet elsewhere in the document.
$netbeansPara.Range.Find.Execute("Netbeans", $true, $false, $false, $false, $false, $true, 1, $false, "Netbeans (8.x)", 2)

# Re-anchor the Toc bookmark so it only spans "Netbeans" again (the Find
# above stretched it to cover the appended " (8.x)" text too).
$bmRange = $d.Range($headingStart, $headingStart + 8)
$d.Bookmarks.Add("_Toc257918838", $bmRange)

# Re-fetch the paragraph's current extent now that it is longer.
$paraEnd = $netbeansPara.Range.End

# Add a collapsed "_GoBack" bookmark right after " (8.x)", i.e. at the very
# end of the paragraph's text (just before its paragraph mark). Inserting a
# bookmark directly at that boundary is unreliable, so place a one-character
# placeholder there first, anchor the bookmark just in front of it, then
# remove the placeholder again -- the bookmark stays put, now correctly
# collapsed at the paragraph's end.
$placeholderPos = $paraEnd - 1
$placeholderRange = $d.Range($placeholderPos, $placeholderPos)
$placeholderRange.InsertAfter("X")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackAnchor = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $goBackAnchor)

$placeholderCharRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderCharRange.Delete()
